$wb = $excel.ActiveWorkbook

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/366fa35f5505d46787f2b8d90384f42e91823f9a/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8e203ab92a766f986091f264eee1e716fbf1978b/e2e/b.md."

# --- Overview sheet: b.md row (row 3) ---
$ws = $wb.Worksheets.Item("Overview")
$ws.Range("E3").Value = "Ready for handoff"
$ws.Range("F3").Value = "Ready for handoff"
$ws.Range("G3").Value = "2016-08-21 22:47:57"

# --- zh-cn sheet: b.md row (row 3) ---
$ws = $wb.Worksheets.Item("zh-cn")
$ws.Range("C3").Value = "Ready for handoff"
$ws.Range("F3").Value = "'False"
$ws.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$ws.Range("H3").Value = "2016-08-21 22:47:52"
$ws.Range("P3").Value = $errorDetail
$ws.Columns.Item(16).ColumnWidth = 40

# --- de-de sheet: b.md row (row 3) ---
$ws = $wb.Worksheets.Item("de-de")
$ws.Range("C3").Value = "Ready for handoff"
$ws.Range("F3").Value = "'False"
$ws.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$ws.Range("H3").Value = "2016-08-21 22:47:57"
$ws.Range("P3").Value = $errorDetail
$ws.Columns.Item(16).ColumnWidth = 40
